# Updated cryptos list on Mon Oct 23 11:19:35 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "30.548.28"
$ws.Cells.Item(2, 5).Value = "  +2.09%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "1.672.55"
$ws.Cells.Item(3, 5).Value = "  +2.45%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  -0.18%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'219.51"
$ws.Cells.Item(5, 5).Value = "  +2.40%  "

# Row 6 - XRP
$ws.Cells.Item(6, 4).Value = "'0.529"
$ws.Cells.Item(6, 5).Value = "  +2.19%  "

# Row 7 - USDC
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 5).Value = "  -0.15%  "

# Row 8 - Solana
$ws.Cells.Item(8, 4).Value = "'29.55"
$ws.Cells.Item(8, 5).Value = "  +3.60%  "

# Row 9 - Cardano
$ws.Cells.Item(9, 5).Value = "  +2.75%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 4).Value = "'0.0640"
$ws.Cells.Item(10, 5).Value = "  +5.35%  "

# Row 11 - TRON
$ws.Cells.Item(11, 4).Value = "'0.0906"
$ws.Cells.Item(11, 5).Value = "  -0.21%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).Value = "1.912.61"
$ws.Cells.Item(12, 5).Value = "  +2.44%  "

# Row 13 - was Polygon, now WrappedEther
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).Value = "1.677.81"
$ws.Cells.Item(13, 5).Value = "  +2.53%  "

# Row 14 - was WrappedEther, now Polygon
$ws.Cells.Item(14, 2).Value = "Polygon"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(14, 4).Value = "'0.613"
$ws.Cells.Item(14, 5).Value = "  +9.05%  "

# Row 15 - Chainlink
$ws.Cells.Item(15, 4).Value = "'10.23"
$ws.Cells.Item(15, 5).Value = "  +10.36%  "

# Row 16 - Polkadot
$ws.Cells.Item(16, 5).Value = "  +4.08%  "

# Row 17 - WrappedBTC
$ws.Cells.Item(17, 4).Value = "30.569.76"
$ws.Cells.Item(17, 5).Value = "  +2.11%  "

# Row 18 - Litecoin
$ws.Cells.Item(18, 4).Value = "'66.36"
$ws.Cells.Item(18, 5).Value = "  +3.71%  "

# Row 19 - BitcoinCash
$ws.Cells.Item(19, 4).Value = "'242.94"
$ws.Cells.Item(19, 5).Value = "  +0.49%  "

# Row 20 - ShibaInu
$ws.Cells.Item(20, 4).Value = "'0.0" + [char]0x2083 + "0722"
$ws.Cells.Item(20, 5).Value = "  +3.10%  "

# Row 21 - Dai
$ws.Cells.Item(21, 5).Value = "  -0.14%  "

# Row 22 - Uniswap
$ws.Cells.Item(22, 5).Value = "  +3.57%  "

# Row 23 - Avalanche
$ws.Cells.Item(23, 5).Value = "  +2.22%  "

# Row 24 - Toncoin
$ws.Cells.Item(24, 4).Value = "'2.17"
$ws.Cells.Item(24, 5).Value = "  +0.10%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'158.39"
$ws.Cells.Item(25, 5).Value = "  +0.24%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'15.85"
$ws.Cells.Item(26, 5).Value = "  +2.23%  "

# Row 27 - Cosmos
$ws.Cells.Item(27, 5).Value = "  +2.61%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +1.09%  "

# Row 29 - BinanceUSD
$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 5).Value = "  -0.18%  "

# Row 30 - Hedera
$ws.Cells.Item(30, 5).Value = "  +1.79%  "

# Row 31 - PancakeSwap
$ws.Cells.Item(31, 4).Value = "'1.15"
$ws.Cells.Item(31, 5).Value = "  +3.05%  "

# Row 32 - Filecoin
$ws.Cells.Item(32, 4).Value = "'3.46"
$ws.Cells.Item(32, 5).Value = "  +2.74%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).Value = "'3.29"
$ws.Cells.Item(33, 5).Value = "  +3.45%  "

# Row 34 - Maker
$ws.Cells.Item(34, 4).Value = "1.489.08"
$ws.Cells.Item(34, 5).Value = "  +4.52%  "

# Row 35 - LidoDAOToken
$ws.Cells.Item(35, 5).Value = "  +7.02%  "

# Row 36 - Aave
$ws.Cells.Item(36, 4).Value = "'84.40"
$ws.Cells.Item(36, 5).Value = "  +11.80%  "

# Row 37 - TrustWalletToken
$ws.Cells.Item(37, 5).Value = "  -0.52%  "

# Row 38 - ImmutableX
$ws.Cells.Item(38, 5).Value = "  +8.26%  "

# Row 39 - VeChain
$ws.Cells.Item(39, 5).Value = "  +5.36%  "

# Row 40 - MXToken
$ws.Cells.Item(40, 5).Value = "  -4.40%  "

# Row 41 - HuobiToken
$ws.Cells.Item(41, 5).Value = "  -0.32%  "

# Row 42 - ARBITRUM
$ws.Cells.Item(42, 5).Value = "  +1.38%  "

# Row 43 - Kaspa
$ws.Cells.Item(43, 5).Value = "  +1.54%  "

# Row 44 - RenderToken
$ws.Cells.Item(44, 5).Value = "  -1.80%  "

# Row 45 - WEMIXToken
$ws.Cells.Item(45, 5).Value = "  -0.20%  "

# Row 46 - PaxDollar
$ws.Cells.Item(46, 4).Value = "'0.999"
$ws.Cells.Item(46, 5).Value = "  -0.16%  "

# Row 47 - FraxShare
$ws.Cells.Item(47, 4).Value = "'5.50"
$ws.Cells.Item(47, 5).Value = "  +2.51%  "

# Row 48 - BitcoinSV
$ws.Cells.Item(48, 4).Value = "'51.06"
$ws.Cells.Item(48, 5).Value = "  -3.23%  "

# Row 49 - RocketPoolETH
$ws.Cells.Item(49, 4).Value = "1.805.69"
$ws.Cells.Item(49, 5).Value = "  +1.77%  "

# Row 50 - Quant
$ws.Cells.Item(50, 4).Value = "'94.91"
$ws.Cells.Item(50, 5).Value = "  +5.01%  "

# Row 51 - BabyDogeCoin
$ws.Cells.Item(51, 5).Value = "  -0.90%  "
